$d = $word.ActiveDocument

# 1. Fix the font name globally (TimesNewToman -> Times New Roman) across all runs
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# 2. Title / author / email replacements (paragraphs 1-3)
$null = $d.Content.Find.Execute('Veiling Revelations: Encryption in the Digital Age', $true, $false, $false, $false, $false, $true, 1, $false, 'Navigating the Labyrinth of Government: A Citizen''s Guide to Civic Engagement', 2)
$null = $d.Content.Find.Execute('Jackson Taylor', $true, $false, $false, $false, $false, $true, 1, $false, 'Amanda Carter', 2)

# Email paragraph (paragraph 3): restrict Find scope to that paragraph to avoid stray matches
$p3 = $d.Paragraphs.Item(3).Range
$null = $p3.Find.Execute('jacktaylor@protonmail', $true, $false, $false, $false, $false, $true, 1, $false, 'amanda', 2)
$p3 = $d.Paragraphs.Item(3).Range
$null = $p3.Find.Execute('com', $true, $false, $false, $false, $false, $true, 1, $false, 'carter@highschool', 2)
# $p3 now equals the replaced "carter@highschool" range; append the two new runs "." and "edu" matching its formatting
$p3.Collapse(0)
$p3.InsertAfter(".")
$p3.Font.Name = "Times New Roman"
$p3.Font.Size = 16
$p3.Font.Color = 0
$p3.Collapse(0)
$p3.InsertAfter("edu")
$p3.Font.Name = "Times New Roman"
$p3.Font.Size = 16
$p3.Font.Color = 0

# 3. Body paragraph replacements (paragraph 5 - essay body)
$null = $d.Content.Find.Execute('In an era where digital footprints shape our lives, the quest for privacy and data security has become paramount', $true, $false, $false, $false, $false, $true, 1, $false, 'In the intricate tapestry of human societies, where laws are woven and leaders emerge, there exists a complex system of governance', 2)
$null = $d.Content.Find.Execute(' Encryption, a process of transforming information into an unreadable format, has emerged as a pivotal tool in safeguarding digital realms, compelling us to delve into its intricacies', $true, $false, $false, $false, $false, $true, 1, $false, ' Government, the entity entrusted with the formidable task of shaping our collective destinies, holds the key to unlocking a harmonious and just society', 2)
$null = $d.Content.Find.Execute(' Encryption''s origins can be traced back to ancient civilizations, with rudimentary techniques employed to protect sensitive communications', $true, $false, $false, $false, $false, $true, 1, $false, ' As high school students, embarking on the cusp of adulthood, it is imperative that we delve into the intricacies of government, deciphering its enigmatic mechanisms and comprehending its profound impact on our lives', 2)
$null = $d.Content.Find.Execute(' As technology evolved, so did encryption methods, culminating in the sophisticated algorithms that underpin modern digital security', $true, $false, $false, $false, $false, $true, 1, $false, ' In this essay, we will embark on a journey through the labyrinthine corridors of government, unraveling its mysteries and illuminating the pathways to civic engagement', 2)
$null = $d.Content.Find.Execute('Encryption has become an indispensable component of internet communication, e-commerce transactions, and countless other digital interactions', $true, $false, $false, $false, $false, $true, 1, $false, 'In our ever-evolving world, where the winds of change blow incessantly, governments stand as the guardians of stability and progress', 2)
$null = $d.Content.Find.Execute(' Its applications extend far beyond the digital sphere, safeguarding sensitive information in industries such as finance, healthcare, and government', $true, $false, $false, $false, $false, $true, 1, $false, ' From the ancient city-states of Greece to the sprawling democracies of today, governments have served as the crucibles of decision-making, balancing the delicate scales of justice and safeguarding the rights and freedoms of citizens', 2)
$null = $d.Content.Find.Execute(' The encryption landscape is dynamic and multifaceted, encompassing various algorithms, key management techniques, and security protocols, each tailored to specific security requirements', $true, $false, $false, $false, $false, $true, 1, $false, ' Yet, these intricate structures are only as effective as the citizens who participate in them and hold them accountable', 2)
$null = $d.Content.Find.Execute(' Understanding the intricacies of encryption is not only essential for safeguarding digital assets but also crucial for navigating the complex regulatory and ethical issues surrounding data protection', $true, $false, $false, $false, $false, $true, 1, $false, ' It is through active engagement that we, as individuals, can shape the course of governance and steer our communities toward a brighter future', 2)
$null = $d.Content.Find.Execute('Encryption serves as a double-edged sword, providing both protection and anonymity', $true, $false, $false, $false, $false, $true, 1, $false, 'As we navigate the complexities of government, it is essential to recognize the diversity of its forms', 2)
$null = $d.Content.Find.Execute(' While it shields data from unauthorized access, it can also conceal illicit activities, making it challenging for law enforcement and intelligence agencies to investigate crimes', $true, $false, $false, $false, $false, $true, 1, $false, ' From the Westminster system, with its monarchy and parliamentary structure, to the presidential system, characterized by the separation of powers, governments across the globe exhibit myriad variations', 2)
$null = $d.Content.Find.Execute(' Balancing the need for privacy with the imperative for public safety has become a contentious debate, pitting the rights of individuals against the collective security of society', $true, $false, $false, $false, $false, $true, 1, $false, ' Understanding these differences is paramount to comprehending the nuances of governance and the unique challenges and opportunities that each system presents', 2)
$null = $d.Content.Find.Execute(' As the digital landscape continues to expand, the significance of encryption will only intensify, demanding ongoing discourse and innovation to address the evolving challenges of data protection', $true, $false, $false, $false, $false, $true, 1, $false, ' It is through comparative analysis and critical reflection that we can cultivate a sophisticated understanding of government, appreciating its complexities while recognizing the common threads that bind all societies together', 2)

# 4. Summary paragraph replacements (paragraph 7)
$null = $d.Content.Find.Execute('Encryption stands as a cornerstone of digital security, protecting sensitive information from unauthorized access', $true, $false, $false, $false, $false, $true, 1, $false, 'In this essay, we explored the labyrinthine world of government, uncovering its intricacies and illuminating the pathways to civic engagement', 2)
$null = $d.Content.Find.Execute(' Its applications span a diverse range of domains, from internet communication to financial transactions', $true, $false, $false, $false, $false, $true, 1, $false, ' We delved into the historical foundations of governance, tracing its evolution from ancient city-states to modern democracies', 2)
$null = $d.Content.Find.Execute(' The encryption landscape is complex and dynamic, ', $true, $false, $false, $false, $false, $true, 1, $false, ' We examined ', 2)
$null = $d.Content.Find.Execute('encompassing a multitude of encryption algorithms, key management techniques, and security protocols, each designed for specific security needs', $true, $false, $false, $false, $false, $true, 1, $false, 'the diverse forms of government, highlighting the Westminster and presidential systems as prominent examples', 2)
$null = $d.Content.Find.Execute(' While encryption safeguards privacy, it also poses challenges in the realm of law enforcement and intelligence gathering', $true, $false, $false, $false, $false, $true, 1, $false, ' Through comparative analysis and critical reflection, we gained a deeper understanding of the challenges and opportunities that each system presents', 2)
$null = $d.Content.Find.Execute(' Striking a balance between individual privacy and public safety remains an ongoing debate', $true, $false, $false, $false, $false, $true, 1, $false, ' As high school students, it is imperative that we embrace our role as active citizens, engaging with government through voting, advocacy, and community involvement', 2)
$null = $d.Content.Find.Execute(' As the digital landscape expands, encryption''s significance will only grow, necessitating continued innovation and discourse to address the evolving challenges posed by data protection', $true, $false, $false, $false, $false, $true, 1, $false, ' By participating in the democratic process, we can shape the course of governance and contribute to the creation of a just and equitable society', 2)

# 5. Add a new empty paragraph at the end of the document (after the Summary paragraph)
$null = $d.Paragraphs.Add()

